# Final  submitted code and Report
#
# Adds an "insertionSort(Integers)" results column, a new 10,000-element
# benchmark row, and a thick outline border around the whole results
# table (A2:D9) with the "# of elements" column right aligned.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room for the new 10,000-element row: insert a row at 3,
#        pushing the existing data rows (old 3-8) down to 4-9. The
#        header row (2) and the blank row above it (1) are untouched. ---
$ws.Rows("3:3").Insert()

# --- 2. Header row: add the new column D header. ---
$ws.Range("D2").Value = "insertionSort(Integers)"

# --- 3. New row 3: the 10,000 element benchmark. ---
$ws.Range("A3").Value = 10000
$ws.Range("B3").Value = "13 msec."
$ws.Range("C3").Value = "28 msec."
$ws.Range("D3").Value = "617 msec."

# --- 4. Existing data rows (now shifted to rows 4-9): fill in the new
#        insertionSort column - anything beyond ~1,000,000 elements
#        took longer than 2 minutes. ---
$ws.Range("D4").Value = "> 2 min"
$ws.Range("D5").Value = "> 2 min"
$ws.Range("D6").Value = "> 2 min"
$ws.Range("D7").Value = "> 2 min"
$ws.Range("D8").Value = "> 2 min"
$ws.Range("D9").Value = "> 2 min"

# --- 5. Column widths for the table. ---
$ws.Columns("A").ColumnWidth = 18.54
$ws.Columns("B").ColumnWidth = 14.26
$ws.Columns("D").ColumnWidth = 19.73

# --- 6. Row height tweak on the 1,000,000-row (matches source report). ---
$ws.Rows(4).RowHeight = 14

# --- 7. Thick outline border around the whole table A2:D9. ---
$tableRange = $ws.Range("A2:D9")
$tableRange.Borders(7).Weight = -4138
$tableRange.Borders(8).Weight = -4138
$tableRange.Borders(9).Weight = -4138
$tableRange.Borders(10).Weight = -4138

# --- 8. Right align the "# of elements" column header + first value. ---
$ws.Range("A2:A3").HorizontalAlignment = -4152

# --- 9. Selection bookkeeping to match the saved view. ---
$ws.Range("A2:D9").Select()
